$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table: row 7 is a new data row, modelled on row 6 (same
# overall shape - some columns such as L/R/S stay blank). Copy row 6's
# formatting down into row 7 first so styles (text / date number formats)
# match exactly, then overwrite with the new row's values.
$ws.Range("A6:AA6").Copy()
$ws.Range("A7:AA7").PasteSpecial(-4122)
$ws.Range("L7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()

$ws.Range("A7").Value = "ABC300"
$ws.Range("B7").Value = "5-"
$ws.Range("C7").Value = "5-"
$ws.Range("D7").Value = "5-"
$ws.Range("E7").Value = "5-"
$ws.Range("F7").Value = "5-"
$ws.Range("G7").Value = "100000002202025"
$ws.Range("H7").Value = "789000002202025 X"
$ws.Range("I7").Value = "ABC300"
$ws.Range("J7").Value = "UNSECURED"
$ws.Range("K7").Value = "01 - Stage 1"
$ws.Range("M7").Value = "5/31/2025"
$ws.Range("N7").Value = "6/20/2058"
$ws.Range("O7").Value = 0.011
$ws.Range("P7").Value = "A3B"
$ws.Range("Q7").Value = "A3B"
$ws.Range("T7").Value = 0.45
$ws.Range("U7").Value = "M_LINEAR"
$ws.Range("V7").Value = 1799989
$ws.Range("W7").Value = "B"
$ws.Range("X7").Value = "5/31/2025"
$ws.Range("Y7").Value = "EUR"
$ws.Range("Z7").Value = 176
$ws.Range("AA7").Value = 0.4

# Update the view: scroll so column H is the left-most visible column and
# select N9 (matches the author's on-screen state when they saved).
$excel.ActiveWindow.ScrollColumn = 8
$ws.Range("N9").Select()
